# Apply timetable room/code corrections across the three timetable sheets
# (Regular_Timetable, PreMid_Timetable, PostMid_Timetable) in the
# sem3_ECE_timetable workbook.

$wb = $excel.ActiveWorkbook

# ---- Regular_Timetable ----
$ws = $wb.Worksheets.Item("Regular_Timetable")
$ws.Range("B2").Value = "MINOR: Generative Ai [C102]"
$ws.Range("E6").Value = "EC263 (Lab) [L107]"
$ws.Range("E7").Value = "EC263 (Lab) [L107]"
$ws.Range("D8").Value = "EC262 (Lab) [L106]"
$ws.Range("D9").Value = "EC262 (Lab) [L106]"
$ws.Range("B10").Value = "MINOR: VLSI [C102]"

# ---- PreMid_Timetable ----
$ws = $wb.Worksheets.Item("PreMid_Timetable")
$ws.Range("B2").Value = "MINOR: Generative Ai [C102]"
$ws.Range("E8").Value = "EC262 (Lab) [L106]"
$ws.Range("E9").Value = "EC262 (Lab) [L106]"
$ws.Range("B10").Value = "MINOR: VLSI [C102]"

# ---- PostMid_Timetable ----
$ws = $wb.Worksheets.Item("PostMid_Timetable")
$ws.Range("B2").Value = "MINOR: Generative Ai [C102]"
$ws.Range("E8").Value = "EC263 (Lab) [L107]"
$ws.Range("E9").Value = "EC263 (Lab) [L107]"
$ws.Range("B10").Value = "MINOR: VLSI [C102]"
